# Add 2022-Q4 data: insert a new quarter sheet right after "总计", shifting
# the existing "2022-Q2" -> "2022-Q4" (with fresh Q4 numbers) and pushing the
# previous "2022-Q2"/"2022-Q1" content into brand-new sheets in their place.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")
$q1 = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert the new 2022-Q4 row at the top of the data
#    and shift the old rows down, preserving their original styles exactly by
#    copying whole rows instead of retyping values.
# ---------------------------------------------------------------------------

# old row 3 (2022-Q1, 0.03) -> row 4
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A4").Value = 2

# old row 2 (2022-Q2, 0.03) -> row 3
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

# row 2 becomes the new 2022-Q4 entry
$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.04

# ---------------------------------------------------------------------------
# 2) Re-create the per-quarter detail sheets in the right order so sheetIds
#    line up as: 总计=1, 2022-Q4=2, 2022-Q2=3, 2022-Q1=4.
#    Delete Q1 first (frees its id), then duplicate Q2 twice (Copy keeps all
#    formatting/styles identical) to rebuild Q2 and Q1 in their new slots.
# ---------------------------------------------------------------------------

$q1.Delete()

# Duplicate "2022-Q2" right after itself -> becomes the (untouched) new Q2.
$q2.Copy($null, $q2)
$newQ2 = $wb.Worksheets.Item("2022-Q2 (2)")

# Duplicate that copy again -> becomes the new Q1 (still holding Q2 data for now).
$newQ2.Copy($null, $newQ2)
$newQ1 = $wb.Worksheets.Item("2022-Q2 (2) (2)")

# Rename in final left-to-right order.
$q2.Name = "2022-Q4"
$newQ2.Name = "2022-Q2"
$newQ1.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 3) Update the values on the new "2022-Q4" sheet (the renamed original Q2).
# ---------------------------------------------------------------------------

$q4 = $q2
$q4.Range("D2").Value = "1.18"
$q4.Range("E2").Value = "92.77"
$q4.Range("F2").Value = "3.25"
$q4.Range("G2").Value = "0.0384"
$q4.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# 4) Restore the 2022-Q1 sheet's own numbers (the copy currently still has
#    the Q2 figures) and its original tab selection.
# ---------------------------------------------------------------------------

$newQ1.Range("D2").Value = "1.35"
$newQ1.Range("E2").Value = "88.71"
$newQ1.Range("F2").Value = "2.55"
$newQ1.Range("G2").Value = "0.0344"
$newQ1.Range("H2").Value = 6
$newQ1.Activate()
